$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell "D2" "42.858.68"

Set-TextCell "D3" "2.282.63"
Set-TextCell "E3" "  +4.79%  "

Set-TextCell "E4" "  +0.06%  "

Set-TextCell "D5" "252.67"
Set-TextCell "E5" "  +0.64%  "

Set-TextCell "E6" "  +4.25%  "

Set-TextCell "D7" "72.84"
Set-TextCell "E7" "  +9.31%  "

Set-TextCell "E8" "  -0.10%  "

Set-TextCell "E9" "  +12.75%  "

Set-TextCell "D10" "38.69"
Set-TextCell "E10" "  +5.99%  "

Set-TextCell "B11" "Dogecoin"
Set-TextCell "C11" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D11" "0.0978"
Set-TextCell "E11" "  +4.33%  "

Set-TextCell "B12" "OKB"
Set-TextCell "C12" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D12" "59.90"
Set-TextCell "E12" "  +1.20%  "

Set-TextCell "D13" "7.37"
Set-TextCell "E13" "  +7.15%  "

Set-TextCell "D14" "0.105"
Set-TextCell "E14" "  +1.43%  "

Set-TextCell "D15" "2.623.21"
Set-TextCell "E15" "  +4.75%  "

Set-TextCell "D16" "14.92"
Set-TextCell "E16" "  +3.93%  "

Set-TextCell "D17" "0.888"
Set-TextCell "E17" "  +4.28%  "

Set-TextCell "D18" "2.288.26"
Set-TextCell "E18" "  +5.40%  "

Set-TextCell "D19" "42.811.97"
Set-TextCell "E19" "  +4.08%  "

Set-TextCell "E20" "  +7.09%  "

Set-TextCell "E21" "  +4.28%  "

Set-TextCell "D22" "73.50"
Set-TextCell "E22" "  +2.47%  "

Set-TextCell "D23" "237.19"
Set-TextCell "E23" "  +2.83%  "

Set-TextCell "D24" "2.13"
Set-TextCell "E24" "  +3.51%  "

Set-TextCell "D25" "3.90"
Set-TextCell "E25" "  +2.12%  "

Set-TextCell "E26" "  +1.41%  "

Set-TextCell "E27" "  -0.10%  "

Set-TextCell "E28" "  +1.69%  "

Set-TextCell "E29" "  -1.03%  "

Set-TextCell "E30" "  +0.35%  "

Set-TextCell "D31" "168.06"
Set-TextCell "E31" "  -0.19%  "

Set-TextCell "D32" "21.06"
Set-TextCell "E32" "  +3.89%  "

Set-TextCell "D33" "6.40"
Set-TextCell "E33" "  +8.79%  "

Set-TextCell "E34" "  +6.73%  "

Set-TextCell "D35" "0.0810"
Set-TextCell "E35" "  +6.89%  "

Set-TextCell "D36" "30.69"
Set-TextCell "E36" "  +24.24%  "

Set-TextCell "E37" "  +4.56%  "

Set-TextCell "D38" "4.76"
Set-TextCell "E38" "  +19.92%  "

Set-TextCell "E39" "  +5.38%  "

Set-TextCell "E40" "  +2.06%  "

Set-TextCell "E41" "  +5.06%  "

Set-TextCell "D42" "13.31"
Set-TextCell "E42" "  +15.80%  "

Set-TextCell "D43" "6.02"
Set-TextCell "E43" "  +9.45%  "

Set-TextCell "D44" "0.212"
Set-TextCell "E44" "  +11.91%  "

Set-TextCell "D45" "9.18"
Set-TextCell "E45" "  +7.66%  "

Set-TextCell "E46" "  -6.77%  "

Set-TextCell "D47" "61.42"
Set-TextCell "E47" "  +0.90%  "

Set-TextCell "D48" "0.103"
Set-TextCell "E48" "  +2.63%  "

Set-TextCell "E49" "  +3.99%  "

Set-TextCell "E50" "  +0.31%  "

Set-TextCell "E51" "  +4.47%  "
